$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new local-extreme row (Laishevsky district) as row 18,
# copying the formatting of the preceding data row (17) first.
$ws.Range("A17:T17").Copy() | Out-Null
$ws.Range("A18:T18").PasteSpecial(-4122) | Out-Null

$ws.Range("A18").Value = 92634000
$ws.Range("B18").Value = "Лаишевский муниципальный район "
$ws.Range("C18").Value = 2020
$ws.Range("D18").Value = 3.549603174603174
$ws.Range("E18").Value = 0.31671285938491328
$ws.Range("F18").Value = 0.67525595756139944
$ws.Range("G18").Value = 0.78245818184798022
$ws.Range("H18").Value = 0.29197724636653227
$ws.Range("I18").Value = 0.26613632168688722
$ws.Range("J18").Value = 0.44619244393439622
$ws.Range("K18").Value = 0.61670443509226291
$ws.Range("L18").Value = 0.33020286625299722
$ws.Range("M18").Value = 0.1389564476656282
$ws.Range("N18").Value = 0.1488110983657692
$ws.Range("O18").Value = 0.018713679303629879
$ws.Range("P18").Value = 0.55618425418864292
$ws.Range("Q18").Value = 0.33667865273731812
$ws.Range("R18").Value = 0.23257069354532661
$ws.Range("S18").Value = 0.44362833249958911
$ws.Range("T18").Value = 1.0587008965338141

# The row that used to be labeled "Котлас" is renamed to its full
# official name "город Котлас".
$ws.Range("B17").Value = "город Котлас"

# Update the sheet's last active-cell selection.
$ws.Range("B24").Select() | Out-Null
